$d = $word.ActiveDocument

# 1) Remove the whole "Modelo Número 1" paragraph entirely (text run,
#    break run and its own pPr all disappear with it).
$d.Paragraphs(2).Range.Delete()

# 2) Capture the plain text of the (still formatted) first paragraph so
#    we can re-type it with no character formatting. Range.Text on a
#    paragraph includes the trailing line-break (chr 11) and paragraph
#    mark (chr 13) characters, so strip those off.
$introText = $d.Paragraphs(1).Range.Text.TrimEnd([char]11, [char]13)

# 3) Prepend a brand-new, unformatted line-break and text run in front
#    of the old (white/underlined) run. Inserting the break first and
#    then the text before it keeps them as two separate runs rather
#    than Word merging them into a single run that mixes text and
#    <w:br/>.
$brPos = $d.Range(0, 0)
$brPos.InsertBefore("`v")
$txtPos = $d.Range(0, 0)
$txtPos.InsertBefore($introText)

# 4) Split the paragraph right between our new (unformatted) content
#    and the old (formatted) content that is still sitting after it,
#    so the old run ends up alone in its own paragraph.
$splitAt = $introText.Length + 1
$splitPos = $d.Range($splitAt, $splitAt)
$splitPos.InsertBefore("`r")

# 5) Delete that now-isolated paragraph, getting rid of the old
#    white/underlined runs for good.
$d.Paragraphs(2).Range.Delete()
